$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.872.07'
$ws.Range("E2").Value = '  +2.70%  '

$ws.Range("D3").Value = '2.609.71'
$ws.Range("E3").Value = '  +1.36%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''579.39'

$ws.Range("D6").Value = '''143.80'
$ws.Range("E6").Value = '  +1.45%  '

$ws.Range("E7").Value = '  -0.26%  '

$ws.Range("E8").Value = '  +0.52%  '

$ws.Range("D9").Value = '2.636.56'
$ws.Range("E9").Value = '  +2.21%  '

$ws.Range("E10").Value = '  -3.06%  '

$ws.Range("E11").Value = '  +2.38%  '

$ws.Range("E12").Value = '  -4.89%  '

$ws.Range("E13").Value = '  +5.77%  '

$ws.Range("D14").Value = '3.079.13'
$ws.Range("E14").Value = '  +1.72%  '

$ws.Range("D15").Value = '60.845.87'
$ws.Range("E15").Value = '  +2.67%  '

$ws.Range("D16").Value = '''23.40'
$ws.Range("E16").Value = '  +1.78%  '

$ws.Range("E17").Value = '  +4.32%  '

$ws.Range("D18").Value = '2.625.37'

$ws.Range("D19").Value = '''11.30'
$ws.Range("E19").Value = '  +9.48%  '

$ws.Range("E20").Value = '  +2.97%  '

$ws.Range("D21").Value = '''350.67'
$ws.Range("E21").Value = '  +3.77%  '

$ws.Range("D22").Value = '''6.93'
$ws.Range("E22").Value = '  +7.27%  '

$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").Value = '''0.518'
$ws.Range("E24").Value = '  +8.36%  '

$ws.Range("D25").Value = '''63.24'
$ws.Range("E25").Value = '  +1.27%  '

$ws.Range("E26").Value = '  -0.20%  '

$ws.Range("E27").Value = '  +0.72%  '

$ws.Range("D28").Value = '''7.91'
$ws.Range("E28").Value = '  +7.31%  '

$ws.Range("D29").Value = '0.0₃0800'
$ws.Range("E29").Value = '  +3.54%  '

$ws.Range("E30").Value = '  +9.05%  '

$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '''6.37'
$ws.Range("E31").Value = '  +3.05%  '

$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D32").Value = '''0.997'
$ws.Range("E32").Value = '  -0.12%  '

$ws.Range("D33").Value = '''162.74'
$ws.Range("E33").Value = '  +2.22%  '

$ws.Range("D34").Value = '''19.59'
$ws.Range("E34").Value = '  +2.82%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '''4.30'
$ws.Range("E35").Value = '  +5.43%  '

$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").Value = '''1.01'
$ws.Range("E36").Value = '  +12.45%  '

$ws.Range("E37").Value = '  +6.52%  '

$ws.Range("D38").Value = '''1.63'
$ws.Range("E38").Value = '  +9.86%  '

$ws.Range("D39").Value = '''37.97'
$ws.Range("E39").Value = '  +1.50%  '

$ws.Range("D40").Value = '''3.90'
$ws.Range("E40").Value = '  +6.24%  '

$ws.Range("D41").Value = '''308.02'
$ws.Range("E41").Value = '  +6.56%  '

$ws.Range("D42").Value = '''0.848'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").Value = '''134.20'
$ws.Range("E43").Value = '  -2.69%  '

$ws.Range("D44").Value = '''20.47'
$ws.Range("E44").Value = '  +9.85%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''0.995'
$ws.Range("E45").Value = '  -0.35%  '

$ws.Range("D46").Value = '''19.94'
$ws.Range("E46").Value = '  +5.69%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '''0.610'
$ws.Range("E47").Value = '  +2.99%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''5.05'
$ws.Range("E48").Value = '  +11.77%  '

$ws.Range("E49").Value = '  +1.23%  '

$ws.Range("D50").Value = '''0.0550'
$ws.Range("E50").Value = '  +3.85%  '

$ws.Range("E51").Value = '  +3.81%  '
